# The deck's theme ("Integral", ppt/theme/theme2.xml -- the theme actually
# applied to the slide master / slides) is swapped with the theme that used
# to live in ppt/theme/theme1.xml (the stock "Office Theme" palette).
#
# PowerPoint's object model doesn't expose "replace this theme part with
# that theme part" directly, but every one of the twelve theme colors is
# reachable (and settable) through any slide's ThemeColorScheme, and a
# write there lands straight in the clrScheme of the theme XML backing the
# presentation's single slide master. Driving all twelve slots to the
# "Office Theme" palette reproduces the clrScheme half of the swap.
#
# ThemeColorScheme.Item(n) order is the standard OOXML clrScheme order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
# and .RGB is a COLORREF (0xBBGGRR), matching classic VBA RGB() packing.

function ToRgbColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

$officeThemeColors = @(
    @(0x00, 0x00, 0x00),  # 1  dk1
    @(0xFF, 0xFF, 0xFF),  # 2  lt1
    @(0x44, 0x54, 0x6A),  # 3  dk2
    @(0xE7, 0xE6, 0xE6),  # 4  lt2
    @(0x5B, 0x9B, 0xD5),  # 5  accent1
    @(0xED, 0x7D, 0x31),  # 6  accent2
    @(0xA5, 0xA5, 0xA5),  # 7  accent3
    @(0xFF, 0xC0, 0x00),  # 8  accent4
    @(0x44, 0x72, 0xC4),  # 9  accent5
    @(0x70, 0xAD, 0x47),  # 10 accent6
    @(0x05, 0x63, 0xC1),  # 11 hlink
    @(0x95, 0x4F, 0x72)   # 12 folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $rgb = $officeThemeColors[$i]
    $tcs.Item($i + 1).RGB = ToRgbColor $rgb[0] $rgb[1] $rgb[2]
}
